$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C width change (match column B's width)
$ws.Columns("C").ColumnWidth = 16.67

# Row 3 updates
$ws.Range("C3").Value = "Selin Company"
$ws.Range("D3").Value = "500.00₺"
$ws.Range("E3").Value = "13-12-2023 18:39:39"

# Row 4 updates
$ws.Range("B4").Value = "Selin Company"
$ws.Range("C4").Value = "Ahmet Company"
$ws.Range("D4").Value = "200.00₺"
$ws.Range("E4").Value = "13-12-2023 18:48:48"

# Delete row 5 entirely (shift cells up)
$ws.Rows("5").Delete()
